$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.03"
$ws.Range("E2").Value = "'-0.50%"

$ws.Range("D3").Value = "'36.19"
$ws.Range("E3").Value = "'-2.51%"

$ws.Range("D4").Value = "'5.039"
$ws.Range("E4").Value = "'-0.19%"

$ws.Range("D5").Value = "'0.07862"
$ws.Range("E5").Value = "'-0.61%"

$ws.Range("D6").Value = "'2.125"
$ws.Range("E6").Value = "'-3.57%"

$ws.Range("D7").Value = "'7.957"
$ws.Range("E7").Value = "'-0.79%"

$ws.Range("D8").Value = "'0.9214"
$ws.Range("E8").Value = "'-0.93%"

$ws.Range("D9").Value = "'0.09531"
$ws.Range("E9").Value = "'-3.92%"

$ws.Range("D10").Value = "'0.1847"
$ws.Range("E10").Value = "'-1.99%"

$ws.Range("D11").Value = "'0.08721"
$ws.Range("E11").Value = "'0.37%"

$ws.Range("D12").Value = "'0.03617"
$ws.Range("E12").Value = "'0.50%"

$ws.Range("D13").Value = "'0.09928"

$ws.Range("D14").Value = "'0.001438"
$ws.Range("E14").Value = "'-2.96%"

$ws.Range("D15").Value = "'0.005705"
$ws.Range("E15").Value = "'1.24%"

$ws.Range("D16").Value = "'3.463"
$ws.Range("E16").Value = "'-0.06%"

$ws.Range("D17").Value = "'4.144"
$ws.Range("E17").Value = "'2.73%"

$ws.Range("E18").Value = "'17.67%"

$ws.Range("E19").Value = "'-1.82%"

$ws.Range("E20").Value = "'0.19%"

$ws.Range("D21").Value = "'5.182"
$ws.Range("E21").Value = "'4.87%"

$ws.Range("D23").Value = "'0.04569"
$ws.Range("E23").Value = "'-0.38%"

$ws.Range("D24").Value = "'0.001235"
$ws.Range("E24").Value = "'-1.23%"

$ws.Range("D25").Value = "'0.004786"
$ws.Range("E25").Value = "'-8.91%"

$ws.Range("D26").Value = "'0.0001303"
$ws.Range("E26").Value = "'-6.81%"

$ws.Range("D27").Value = "'0.0004758"
$ws.Range("E27").Value = "'75.39%"

$ws.Range("D39").Value = "'0.01849"
$ws.Range("E39").Value = "'0.75%"

$ws.Range("D40").Value = "'0.04708"
$ws.Range("E40").Value = "'-1.62%"

$ws.Range("D41").Value = "'0.007787"
$ws.Range("E41").Value = "'-2.54%"

$ws.Range("D42").Value = "'0.1385"
$ws.Range("E42").Value = "'-2.15%"

$ws.Range("D43").Value = "'0.007736"
$ws.Range("E43").Value = "'2.33%"

$ws.Range("D44").Value = "'0.002205"
$ws.Range("E44").Value = "'4.44%"

$ws.Range("D45").Value = "'0.01118"
$ws.Range("E45").Value = "'7.25%"

$ws.Range("D46").Value = "'0.00006390"
$ws.Range("E46").Value = "'1.24%"

$ws.Range("E47").Value = "'0.35%"

$ws.Range("E48").Value = "'0.35%"

$ws.Range("D49").Value = "'52.13"
$ws.Range("E49").Value = "'45.87%"

$ws.Range("D50").Value = "'0.001904"
$ws.Range("E50").Value = "'-29.10%"

$ws.Range("D51").Value = "'0.00002104"
$ws.Range("E51").Value = "'0.35%"
